$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain text in the source sheet (e.g. "26.281.82"),
# so force Text format on the specific cells we touch before writing them -
# otherwise Excel would auto-parse values like "1.007" as a number and drop
# formatting such as trailing zeros (e.g. "0.06370" -> 0.0637).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.281.82'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.668.24'
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.93'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5289'
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  +0.68%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06370'
$ws.Range("E9").Value = '  +0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.93'
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07833'
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.671.53'
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.897.54'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5602'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8103'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.309.25'
$ws.Range("E18").Value = '  +0.24%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '200.15'
$ws.Range("E21").Value = '  +4.21%  '
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.064'
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.46'
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1218'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.239'
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.18'
$ws.Range("E28").Value = '  +0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.529'
$ws.Range("E29").Value = '  +3.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05896'
$ws.Range("E30").Value = '  +0.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.283'
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.520'
$ws.Range("E32").Value = '  -1.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.336'
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.599'
$ws.Range("E34").Value = '  -0.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9634'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.815'
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5809'
$ws.Range("E38").Value = '  +0.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01610'
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.987'
$ws.Range("E40").Value = '  +1.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.078.19'
$ws.Range("E41").Value = '  +3.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8580'
$ws.Range("E42").Value = '  +0.72%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.90'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.808.09'
$ws.Range("E45").Value = '  +0.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '58.53'
$ws.Range("E46").Value = '  +2.58%  '
$ws.Range("E47").Value = '  +0.10%  '
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.059'
$ws.Range("E49").Value = '  +1.16%  '
$ws.Range("E50").Value = '  -4.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05145'
$ws.Range("E51").Value = '  -0.34%  '

Write-Host "Updated cryptos list"
